$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "42.129.08"
$ws.Range("E2").Value = "  +2.64%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.212.35"
$ws.Range("E3").Value = "  +1.51%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.06%  "

# Row 5 - BNB
$ws.Range("D5").Value = "252.21"
$ws.Range("E5").Value = "  +0.82%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.43%  "

# Row 7 - Solana
$ws.Range("D7").Value = "67.86"
$ws.Range("E7").Value = "  +0.77%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.11%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.617"
$ws.Range("E9").Value = "  +8.42%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "38.57"
$ws.Range("E10").Value = "  +4.01%  "

# Row 11 - OKB
$ws.Range("D11").Value = "59.42"
$ws.Range("E11").Value = "  +1.93%  "

# Row 12 - Dogecoin
$ws.Range("D12").Value = "0.0939"
$ws.Range("E12").Value = "  +1.32%  "

# Row 13 - Polkadot
$ws.Range("E13").Value = "  +0.82%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +0.08%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.543.87"
$ws.Range("E15").Value = "  +1.70%  "

# Row 16 - Polygon
$ws.Range("D16").Value = "0.871"
$ws.Range("E16").Value = "  +1.45%  "

# Row 17 - Chainlink
$ws.Range("D17").Value = "14.53"
$ws.Range("E17").Value = "  +0.96%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.206.97"
$ws.Range("E18").Value = "  +1.18%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "41.984.73"
$ws.Range("E19").Value = "  +2.48%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0962"
$ws.Range("E20").Value = "  +2.08%  "

# Row 21/22 swap: Litecoin <-> Uniswap
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "6.15"
$ws.Range("E21").Value = "  -0.45%  "

$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").Value = "72.34"
$ws.Range("E22").Value = "  +1.14%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "231.67"
$ws.Range("E23").Value = "  +0.38%  "

# Row 24 - ImmutableX
$ws.Range("E24").Value = "  -2.37%  "

# Row 25 - WEMIXToken
$ws.Range("D25").Value = "3.87"
$ws.Range("E25").Value = "  +1.89%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.11%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "11.22"
$ws.Range("E27").Value = "  -4.21%  "

# Row 28 - PancakeSwap
$ws.Range("E28").Value = "  -2.43%  "

# Row 29 - LEO
$ws.Range("E29").Value = "  -1.08%  "

# Row 30 - Toncoin
$ws.Range("D30").Value = "2.20"
$ws.Range("E30").Value = "  +0.79%  "

# Row 31 - Monero
$ws.Range("D31").Value = "166.87"
$ws.Range("E31").Value = "  -1.70%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "20.43"
$ws.Range("E32").Value = "  -0.13%  "

# Row 33 - Kaspa
$ws.Range("E33").Value = "  +3.84%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("D34").Value = "5.86"
$ws.Range("E34").Value = "  +8.58%  "

# Row 35 - Hedera
$ws.Range("D35").Value = "0.0781"
$ws.Range("E35").Value = "  +8.40%  "

# Row 36 - Stellar
$ws.Range("E36").Value = "  -0.03%  "

# Row 37/38 swap: Filecoin <-> InjectiveProtocol
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").Value = "25.93"
$ws.Range("E37").Value = "  +2.31%  "

$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").Value = "4.59"
$ws.Range("E38").Value = "  +1.05%  "

# Row 39 - RenderToken
$ws.Range("D39").Value = "4.10"
$ws.Range("E39").Value = "  +2.90%  "

# Row 40 - VeChain
$ws.Range("D40").Value = "0.0312"
$ws.Range("E40").Value = "  +5.78%  "

# Row 41 - LidoDAOToken
$ws.Range("E41").Value = "  +1.03%  "

# Row 42/44 swap: FTXToken <-> Celestia (Row 43 THORChain stays in place)
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "12.08"
$ws.Range("E42").Value = "  +0.37%  "

# Row 43 - THORChain
$ws.Range("D43").Value = "5.65"
$ws.Range("E43").Value = "  +0.47%  "

$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").Value = "5.13"
$ws.Range("E44").Value = "  +7.27%  "

# Row 45 - MultiversX
$ws.Range("D45").Value = "61.46"
$ws.Range("E45").Value = "  -4.11%  "

# Row 46 - Algorand
$ws.Range("E46").Value = "  -1.84%  "

# Row 47 - FraxShare
$ws.Range("E47").Value = "  -0.10%  "

# Row 48 - Cronos
$ws.Range("D48").Value = "0.0996"
$ws.Range("E48").Value = "  -1.77%  "

# Row 49 - BinanceUSD
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  -0.22%  "

# Row 50 - ARBITRUM
$ws.Range("E50").Value = "  +1.70%  "

# Row 51 - HuobiToken
$ws.Range("E51").Value = "  +4.37%  "
